$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 28549
$ws.Range("J68").Value = 28549
$ws.Range("L68").Value = 28549
$ws.Range("N68").Value = -30047
$ws.Range("H69").Value = 4510.5
$ws.Range("I69").Value = 3341.5
$ws.Range("J69").Value = 5095
$ws.Range("K69").Value = 10024.5
$ws.Range("L69").Value = 15285
$ws.Range("M69").Value = -9150.5
$ws.Range("N69").Value = -17033
$ws.Range("H70").Value = 1778.6666
$ws.Range("I70").Value = 997.4167
$ws.Range("J70").Value = 2820.3333
$ws.Range("K70").Value = 2992.2501
$ws.Range("L70").Value = 8460.999899999999
$ws.Range("M70").Value = -2722.2501
$ws.Range("N70").Value = -9000.999899999999
$ws.Range("H71").Value = 28549
$ws.Range("J71").Value = 28549
$ws.Range("L71").Value = 85647
$ws.Range("N71").Value = -93135
$ws.Range("H72").Value = 4510.5
$ws.Range("I72").Value = 3341.5
$ws.Range("J72").Value = 5095
$ws.Range("K72").Value = 30073.5
$ws.Range("L72").Value = 45855
$ws.Range("M72").Value = -25705.5
$ws.Range("N72").Value = -54591
$ws.Range("H73").Value = 1778.6666
$ws.Range("I73").Value = 997.4167
$ws.Range("J73").Value = 2820.3333
$ws.Range("K73").Value = 2992.2501
$ws.Range("L73").Value = 8460.999899999999
$ws.Range("M73").Value = -2056.2501
$ws.Range("N73").Value = -10332.9999
$ws.Range("H98").Value = 3737.1765
$ws.Range("I98").Value = 4310.154
$ws.Range("K98").Value = 4310.154
$ws.Range("M98").Value = -2812.154
$ws.Range("H122").Value = 3737.1765
$ws.Range("I122").Value = 4310.154
$ws.Range("K122").Value = 12930.462
$ws.Range("M122").Value = -10480.462
$ws.Range("H132").Value = 8003394.5
$ws.Range("I132").Value = 10529791
$ws.Range("J132").Value = 3141
$ws.Range("K132").Value = 31589373
$ws.Range("L132").Value = 9423
$ws.Range("M132").Value = -31586843
$ws.Range("N132").Value = -14483

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8930102
$ws.Range("I2").Value = 17858792
$ws.Range("J2").Value = 1412.8572
$ws.Range("K2").Value = 17858792
$ws.Range("L2").Value = 1412.8572
$ws.Range("M2").Value = -17858679
$ws.Range("N2").Value = -1638.8572
$ws.Range("H45").Value = 1235.0769
$ws.Range("I45").Value = 1113.5652
$ws.Range("J45").Value = 2166.6667
$ws.Range("K45").Value = 1113.5652
$ws.Range("L45").Value = 2166.6667
$ws.Range("M45").Value = -736.5652
$ws.Range("N45").Value = -2920.6667
$ws.Range("H116").Value = 8930102
$ws.Range("I116").Value = 17858792
$ws.Range("J116").Value = 1412.8572
$ws.Range("K116").Value = 17858792
$ws.Range("L116").Value = 1412.8572
$ws.Range("M116").Value = -17856498
$ws.Range("N116").Value = -6000.8572
$ws.Range("H122").Value = 1416.45
$ws.Range("I122").Value = 1076.375
$ws.Range("K122").Value = 3229.125
$ws.Range("M122").Value = -779.125
$ws.Range("H132").Value = 1935.0779
$ws.Range("I132").Value = 1805.6377
$ws.Range("K132").Value = 5416.9131
$ws.Range("M132").Value = -2886.9131
$ws.Range("H134").Value = 30912.5
$ws.Range("J134").Value = 30912.5
$ws.Range("L134").Value = 30912.5
$ws.Range("N134").Value = -41052.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8930102
$ws.Range("I3").Value = 17858792
$ws.Range("J3").Value = 1412.8572
$ws.Range("K3").Value = 17858792
$ws.Range("L3").Value = 1412.8572
$ws.Range("M3").Value = -17858678
$ws.Range("N3").Value = -1640.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10565
$ws.Range("I34").Value = 331.25
$ws.Range("J34").Value = 51500
$ws.Range("K34").Value = 993.75
$ws.Range("L34").Value = 154500
$ws.Range("M34").Value = -909.75
$ws.Range("N34").Value = -154668
$ws.Range("H55").Value = 1345.4546
$ws.Range("I55").Value = 640
$ws.Range("J55").Value = 1933.3334
$ws.Range("K55").Value = 1920
$ws.Range("L55").Value = 5800.0002
$ws.Range("M55").Value = -1743
$ws.Range("N55").Value = -6154.0002
$ws.Range("H93").Value = 2877.6667
$ws.Range("H94").Value = 3720.606
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3831.6128
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 11494.8384
$ws.Range("M94").Value = -5324
$ws.Range("N94").Value = -12846.8384
$ws.Range("H97").Value = 2100
$ws.Range("I97").Value = 500
$ws.Range("J97").Value = 2633.3333
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 7899.999899999999
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -8891.999899999999
$ws.Range("H131").Value = 1373.1538
$ws.Range("J131").Value = 1087.5581
$ws.Range("L131").Value = 3262.6743
$ws.Range("N131").Value = -13342.6743

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31872
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99360
$ws.Range("H93").Value = 20422.562
$ws.Range("J93").Value = 20422.562
$ws.Range("L93").Value = 20422.562
$ws.Range("N93").Value = -24166.562
$ws.Range("H122").Value = 3824.3462
$ws.Range("I122").Value = 3809.625
$ws.Range("J122").Value = 3830.889
$ws.Range("K122").Value = 11428.875
$ws.Range("L122").Value = 11492.667
$ws.Range("M122").Value = -8978.875
$ws.Range("N122").Value = -16392.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1533.9231
$ws.Range("I7").Value = 1198
$ws.Range("J7").Value = 1869.8462
$ws.Range("K7").Value = 1198
$ws.Range("L7").Value = 1869.8462
$ws.Range("M7").Value = -1086
$ws.Range("N7").Value = -2093.8462
$ws.Range("H46").Value = 1077.1233
$ws.Range("I46").Value = 948.31665
$ws.Range("J46").Value = 1671.6154
$ws.Range("K46").Value = 948.31665
$ws.Range("L46").Value = 1671.6154
$ws.Range("M46").Value = -760.31665
$ws.Range("N46").Value = -2047.6154
$ws.Range("H126").Value = 1533.9231
$ws.Range("I126").Value = 1198
$ws.Range("J126").Value = 1869.8462
$ws.Range("K126").Value = 3594
$ws.Range("L126").Value = 5609.5386
$ws.Range("M126").Value = -1124
$ws.Range("N126").Value = -10549.5386
$ws.Range("H132").Value = 6187.241
$ws.Range("I132").Value = 1961.6666
$ws.Range("J132").Value = 10714.643
$ws.Range("K132").Value = 5884.9998
$ws.Range("L132").Value = 32143.929
$ws.Range("M132").Value = -3354.9998
$ws.Range("N132").Value = -37203.929
$ws.Range("H136").Value = 2819.0625
$ws.Range("I136").Value = 2930.5
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 8791.5
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").Value = -6241.5
$ws.Range("N136").Value = -12999.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 34032.5
$ws.Range("J75").Value = 34032.5
$ws.Range("L75").Value = 34032.5
$ws.Range("N75").Value = -35904.5
$ws.Range("H78").Value = 34032.5
$ws.Range("J78").Value = 34032.5
$ws.Range("L78").Value = 102097.5
$ws.Range("N78").Value = -111457.5
$ws.Range("H80").Value = 24333.334
$ws.Range("J80").Value = 24333.334
$ws.Range("L80").Value = 24333.334
$ws.Range("N80").Value = -26329.334
$ws.Range("H83").Value = 24333.334
$ws.Range("J83").Value = 24333.334
$ws.Range("L83").Value = 73000.00199999999
$ws.Range("N83").Value = -82984.00199999999
$ws.Range("H126").Value = 2034.6
$ws.Range("I126").Value = 1767.091
$ws.Range("J126").Value = 2487.3076
$ws.Range("K126").Value = 5301.272999999999
$ws.Range("L126").Value = 7461.9228
$ws.Range("M126").Value = -2831.272999999999
$ws.Range("N126").Value = -12401.9228
